# "Add template support on whatsapp"
# Replace the sample Phone/name/charge sheet with Phone/financial_year/final_date
# WhatsApp template data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (row 1): keep "Phone", swap "name" -> "financial_year" and
# "charge" -> "final_date"
$ws.Range("B1").Value = "financial_year"
$ws.Range("C1").Value = "final_date"

# Data rows: both contacts now carry the same financial year / final date
# template values instead of per-person name/charge.
$ws.Range("B2").Value = "2025-26"
$ws.Range("C2").Value = "31.03.2025"

$ws.Range("B3").Value = "2025-26"
$ws.Range("C3").Value = "31.03.2025"

# Column B grew slightly wider to fit "financial_year"/"2025-26".
$ws.Columns.Item(2).ColumnWidth = 11.6
